$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 6 header ("Légende"/key) row gets a custom row height
$ws.Rows.Item(6).RowHeight = 22.5

# Task "Classe Salle" (row 17) starts: moved from "A faire" to "En cours",
# assigned to "Tous", with a start date set.
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = "X"
$ws.Range("E17").Value = "Tous"
$ws.Range("F17").Value = Get-Date -Year 2016 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0

# Update the last-active-cell selection on the sheet
$ws.Range("J11").Select()
